$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 22 (pushes the signature block
#    rows 26/27 down to 27/28). The new row inherits row 21's
#    current ("closing"/last-row) border style, so afterwards we
#    demote row 21 back to the regular "middle row" style (copied
#    from row 20) and leave row 22 with the closing style.
# ------------------------------------------------------------------
$ws.Rows.Item(22).Insert()
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Re-sort the "Periodo Mora" column into ascending order
#    (2502 .. 2507 across the existing rows) and populate the new
#    row 22 with the extra period 2508.
# ------------------------------------------------------------------
$ws.Range("E16").Value = "2502"
$ws.Range("E17").Value = "2503"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2505"
$ws.Range("E20").Value = "2506"
$ws.Range("E21").Value = "2507"

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047454900"
$ws.Range("D22").Value = "MICHELLE SANTOS PABA"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 120000
$ws.Range("G22").Value = 3000000

# ------------------------------------------------------------------
# 3. Update the summary fields: one more period (7 instead of 6)
#    means the total "Valor Mora" grows by another 120000.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 840000
$ws.Range("F13").Value = 7

# ------------------------------------------------------------------
# 4. Fix up the signature block so the underline row (row 27) is
#    consistent on both sides, and the labels sit on row 28.
# ------------------------------------------------------------------
$ws.Range("B27").Value = "___________________________________"
$ws.Range("H27").Value = "___________________________________"

$ws.Range("B28").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H28").Value = "FIRMA DEL REPRESENTANTE LEGAL"
